$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of trade data (row 5)
$ws.Range("A5").Value = 9959.23
$ws.Range("B5").Value = 10035.5
$ws.Range("C5").Value = 109.08
$ws.Range("D5").Value = 108.25
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = -0.76
$ws.Range("G5").Value = 42612.67454861111
$ws.Range("H5").Value = $false
